$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: reorder "Recorded By" email list ---
$ws.Range("G2").Value = "gehanadel@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, System"

# --- Row 3: reorder "Recorded By" email list ---
$ws.Range("G3").Value = "Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, System, eman.tantawi@med.asu.edu.eg"

# --- Row 4: reorder "Recorded By" email list ---
$ws.Range("G4").Value = "gehanadel@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"

# --- Row 5: reorder "Recorded By" email list ---
$ws.Range("G5").Value = "eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"

# --- Row 6: reorder "Recorded By" email list, recorded-sessions count 23 -> 24 ---
$ws.Range("G6").Value = "alshimaa.atef@med.asu.edu.egm, majorelle.magdy@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"
$ws.Range("L6").Value = 24

# --- Row 7: reorder "Recorded By" email list, missing-sessions count 3 -> 2 ---
$ws.Range("G7").Value = "NadaMohamed@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg"
$ws.Range("L7").Value = 2

# --- Row 9 / Row 10: updated percentage text (must stay plain text, not be
#     auto-converted to a numeric percentage by the COM value-setter) ---
$helper = $ws.Range("Z1")
$helper.NumberFormat = "@"

$helper.Value = "82.8%"
$helper.Copy()
$ws.Range("L9").PasteSpecial(-4163)

$helper.Value = "27.4%"
$helper.Copy()
$ws.Range("L10").PasteSpecial(-4163)

$helper.Clear()

# --- Row 12: reorder "Recorded By" email list ---
$ws.Range("G12").Value = "amira.m.ibrahim@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, dina.adel@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg"

# --- Row 15: Group Statistics mirror of the ANATOMY class-statistics block ---
$ws.Range("O15").Value = 24
$ws.Range("P15").Value = 2

$helper = $ws.Range("Z1")
$helper.NumberFormat = "@"

$helper.Value = "82.8%"
$helper.Copy()
$ws.Range("R15").PasteSpecial(-4163)

$helper.Value = "27.4%"
$helper.Copy()
$ws.Range("S15").PasteSpecial(-4163)

$helper.Clear()

# --- Row 17: reorder "Recorded By" email list ---
$ws.Range("G17").Value = "esraa.sami@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"

# --- Row 18: session got recorded - copy the "Recorded" (green) formatting
#     from row 2 onto row 18, then fill in the new values ---
$ws.Range("A2:I2").Copy()
$ws.Range("A18:I18").PasteSpecial(-4122)

$ws.Cells.Item(18, 7).Value = "afnan.fares@med.asu.edu.eg"
$ws.Cells.Item(18, 8).Value = "83/251"
$ws.Cells.Item(18, 9).Value = "Recorded"

# --- Row 24: reorder "Recorded By" email list ---
$ws.Range("G24").Value = "Sarah.Mahdy@med.asu.edu.eg, youstina.gamil@med.asu.edu.eg"

# --- Row 30: reorder "Recorded By" email list ---
$ws.Range("G30").Value = "shorokmohamed@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg"
